$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '56.800.82'
$ws.Range("E2").Value = '  +0.55%  '

# Row 3
$ws.Range("D3").Value = '3.030.67'
$ws.Range("E3").Value = '  +2.51%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '510.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.89%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.80%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("E8").Value = '  +1.85%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.14'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.12%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.108'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.33%  '

# Row 11
$ws.Range("E11").Value = '  +5.41%  '

# Row 12
$ws.Range("D12").Value = '3.550.53'
$ws.Range("E12").Value = '  +2.47%  '

# Row 13
$ws.Range("E13").Value = '  +0.72%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.80%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000163'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.97%  '

# Row 16
$ws.Range("D16").Value = '56.807.59'
$ws.Range("E16").Value = '  +0.37%  '

# Row 17
$ws.Range("D17").Value = '3.034.92'
$ws.Range("E17").Value = '  +2.46%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.21%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.92%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.24%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '334.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.58%  '

# Row 22
$ws.Range("E22").Value = '  +0.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.502'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.63%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.74'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.51%  '

# Row 25
$ws.Range("D25").Value = '3.158.92'
$ws.Range("E25").Value = '  +2.58%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.166'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.07%  '

# Row 27
$ws.Range("E27").Value = '  -0.14%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0932'
$ws.Range("E28").Value = '  +8.32%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.42%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.65%  '

# Row 31
$ws.Range("E31").Value = '  +2.82%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.45'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.77%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.99%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '153.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.42%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.75%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '27.13'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +13.87%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.83'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.70%  '

# Row 38
$ws.Range("E38").Value = '  +2.23%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0664'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.05%  '

# Row 40
$ws.Range("E40").Value = '  +2.87%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.59'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.78%  '

# Row 42
$ws.Range("E42").Value = '  +0.00%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.77%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.662'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.59%  '

# Row 45
$ws.Range("D45").Value = '2.212.77'
$ws.Range("E45").Value = '  +2.98%  '

# Row 46
$ws.Range("E46").Value = '  +0.74%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0244'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.95%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.939'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.02%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.88'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.93%  '

# Row 50
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.11%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0857'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.24%  '
